$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.452.01'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.26%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.849.48'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.34%  '

$ws.Range("E4").Value = '  +0.20%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.97'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.84%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6273'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.63%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.002'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.17%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07673'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.82%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2916'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.45%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.83'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.73%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07753'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.56%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.852.61'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.68%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.033'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.76%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6814'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.42%  '

$ws.Range("E15").Value = '  +3.56%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '83.50'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.66%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.173'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.11%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.446.56'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.12%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '228.62'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.29%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.39'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.17%  '

$ws.Range("E21").Value = '  +0.11%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.422'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.43%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.002'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.11%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '158.11'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.74%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1373'
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.405'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.46%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.70'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.74%  '

$ws.Range("E28").Value = '  +6.14%  '

$ws.Range("E29").Value = '  +0.61%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05651'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.55%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.120'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.43%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.036'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.37%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.843'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.35%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.162'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.47%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7082'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.52%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.594'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.24%  '

$ws.Range("B37").Value = 'Maker'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.226.76'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.29%  '

$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.761'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.15%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01789'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.87%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.548'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.80%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9026'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.30%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.002'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.25%  '

$ws.Range("B43").Value = 'RocketPoolETH'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.009.28'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.90%  '

$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.79'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.24%  '

$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '66.01'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.53%  '

$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000122'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.19%  '

$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.164'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.40%  '

$ws.Range("B48").Value = 'TheSandbox'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4015'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.32%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1155'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.26%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.014'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.21%  '

$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.672'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.24%  '
